$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 503
$ws1.Range("F7").Value = 7551
$ws1.Range("F11").Value = 5478
$ws1.Range("F15").Value = 7551
$ws1.Range("F16").Value = 8898
$ws1.Range("F18").Value = 888
$ws1.Range("F19").Value = 4422
$ws1.Range("F21").Value = 215
$ws1.Range("F26").Value = 106
$ws1.Range("F27").Value = 1650
$ws1.Range("F28").Value = 706
$ws1.Range("F29").Value = 904
$ws1.Range("F31").Value = 1863
$ws1.Range("F33").Value = 2272
$ws1.Range("F36").Value = 1448
$ws1.Range("F41").Value = 4070
$ws1.Range("F45").Value = 507

# Sheet "全部类型" (All types) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 503
$ws4.Range("F11").Value = 5478
$ws4.Range("F13").Value = 7551
$ws4.Range("F16").Value = 888
$ws4.Range("F17").Value = 4422
$ws4.Range("F19").Value = 215
$ws4.Range("F26").Value = 106
$ws4.Range("F27").Value = 1650
$ws4.Range("F28").Value = 1863
$ws4.Range("F30").Value = 2272
$ws4.Range("F41").Value = 4070
$ws4.Range("F46").Value = 507
